$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded for "Ají" (Inferno / Primera,
# Región de Arica y Parinacota, $/caja 12 kilos) dated 2022-01-28 (serial
# 44589). It belongs chronologically before the existing row 114 entry, so
# insert a fresh row there; Excel shifts every row from 114 down through
# 206 down by one (old row 206 becomes the new row 207).
$ws.Rows(114).Insert()

# Populate the newly inserted row 114 with the new observation. The
# descriptive/meta columns mirror the entry that used to sit at row 114
# (now shifted to row 115); only the date, volume, prices and $/Kg differ.
$ws.Range("A114").Value = 4
$ws.Range("B114").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C114").Value = "Los Lagos"
$ws.Range("D114").Value = 44589
$ws.Range("E114").Value = 10
$ws.Range("F114").Value = 100112021
$ws.Range("G114").Value = "Ají"
$ws.Range("H114").Value = "Inferno"
$ws.Range("I114").Value = "Primera"
$ws.Range("J114").Value = 160
$ws.Range("K114").Value = 17000
$ws.Range("L114").Value = 18000
$ws.Range("M114").Value = 17500
$ws.Range("N114").Value = "$/caja 12 kilos"
$ws.Range("O114").Value = "Región de Arica y Parinacota"
$ws.Range("P114").Value = 1458
$ws.Range("Q114").Value = 12
$ws.Range("R114").Value = "Hortaliza"
